$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.692.37'
$ws.Range("E2").Value = '  -4.78%  '

$ws.Range("D3").Value = '1.718.98'
$ws.Range("E3").Value = '  -5.32%  '

$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  -0.33%  '

$ws.Range("D5").Value = '''223.99'
$ws.Range("E5").Value = '  -3.79%  '

$ws.Range("D6").Value = '''0.5667'
$ws.Range("E6").Value = '  -2.90%  '

$ws.Range("D7").Value = '''1.004'
$ws.Range("E7").Value = '  -0.34%  '

$ws.Range("D8").Value = '''0.2691'
$ws.Range("E8").Value = '  -0.12%  '

$ws.Range("D9").Value = '''22.70'
$ws.Range("E9").Value = '  +0.27%  '

$ws.Range("D10").Value = '''0.06536'
$ws.Range("E10").Value = '  -2.60%  '

$ws.Range("D11").Value = '''0.07509'
$ws.Range("E11").Value = '  -0.01%  '

$ws.Range("D12").Value = '1.728.81'
$ws.Range("E12").Value = '  -4.66%  '

$ws.Range("D13").Value = '''4.635'
$ws.Range("E13").Value = '  +0.87%  '

$ws.Range("D14").Value = '''0.5927'
$ws.Range("E14").Value = '  -3.15%  '

$ws.Range("D15").Value = '1.956.24'
$ws.Range("E15").Value = '  -4.66%  '

$ws.Range("D16").Value = '''73.56'
$ws.Range("E16").Value = '  -0.76%  '

$ws.Range("D17").Value = '''0.000008540'
$ws.Range("E17").Value = '  -8.48%  '

$ws.Range("D18").Value = '27.679.86'
$ws.Range("E18").Value = '  -4.10%  '

$ws.Range("D19").Value = '''5.244'
$ws.Range("E19").Value = '  -2.46%  '

$ws.Range("D20").Value = '''1.003'
$ws.Range("E20").Value = '  -0.38%  '

$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").Value = '''11.16'
$ws.Range("E21").Value = '  -0.73%  '

$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").Value = '''199.90'
$ws.Range("E22").Value = '  -2.69%  '

$ws.Range("D23").Value = '''6.494'
$ws.Range("E23").Value = '  -2.78%  '

$ws.Range("D24").Value = '''1.005'
$ws.Range("E24").Value = '  -0.52%  '

$ws.Range("D25").Value = '''149.05'
$ws.Range("E25").Value = '  -3.29%  '

$ws.Range("D26").Value = '''7.924'
$ws.Range("E26").Value = '  +2.89%  '

$ws.Range("D27").Value = '''0.1208'
$ws.Range("E27").Value = '  -3.21%  '

$ws.Range("D28").Value = '''15.96'

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '''1.376'
$ws.Range("E29").Value = '  -2.24%  '

$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '''0.06115'
$ws.Range("E30").Value = '  -2.35%  '

$ws.Range("D31").Value = '''1.379'
$ws.Range("E31").Value = '  -3.74%  '

$ws.Range("D32").Value = '''3.672'
$ws.Range("E32").Value = '  -0.17%  '

$ws.Range("D33").Value = '''3.671'
$ws.Range("E33").Value = '  +1.22%  '

$ws.Range("D34").Value = '''1.658'
$ws.Range("E34").Value = '  -1.07%  '

$ws.Range("D35").Value = '''1.022'
$ws.Range("E35").Value = '  -2.02%  '

$ws.Range("D36").Value = '''0.6382'
$ws.Range("E36").Value = '  +1.93%  '

$ws.Range("D37").Value = '''2.424'
$ws.Range("E37").Value = '  -4.41%  '

$ws.Range("D38").Value = '''2.674'
$ws.Range("E38").Value = '  -2.58%  '

$ws.Range("D39").Value = '''0.01647'
$ws.Range("E39").Value = '  -2.70%  '

$ws.Range("D40").Value = '1.110.67'
$ws.Range("E40").Value = '  -1.04%  '

$ws.Range("D41").Value = '''6.128'
$ws.Range("E41").Value = '  -3.49%  '

$ws.Range("D42").Value = '''0.8712'
$ws.Range("E42").Value = '  +1.52%  '

$ws.Range("E43").Value = '  -0.18%  '

$ws.Range("D44").Value = '''99.14'
$ws.Range("E44").Value = '  -0.32%  '

$ws.Range("D45").Value = '1.869.77'
$ws.Range("E45").Value = '  -5.05%  '

$ws.Range("D46").Value = '''58.52'
$ws.Range("E46").Value = '  -1.93%  '

$ws.Range("D47").Value = '''0.00000000109'
$ws.Range("E47").Value = '  -3.53%  '

$ws.Range("D48").Value = '''1.542'
$ws.Range("E48").Value = '  -1.12%  '

$ws.Range("D49").Value = '''8.192'
$ws.Range("E49").Value = '  -0.40%  '

$ws.Range("D50").Value = '''0.05362'
$ws.Range("E50").Value = '  -2.39%  '

$ws.Range("D51").Value = '''0.4407'
$ws.Range("E51").Value = '  -2.96%  '
